$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix comma used instead of period as a separator in scraped contractor names
$ws.Range("E30").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F30").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E54").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("F54").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("E61").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E62").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F62").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E63").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E80").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F80").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"

# Fix floating point "Importe" values scraped with comma decimal separator and
# dot thousands separator -> plain dot-decimal text (leading apostrophe keeps it text)
$ws.Range("H2").Value = "'520.00"
$ws.Range("H3").Value = "'480.00"
$ws.Range("H4").Value = "'30.00"
$ws.Range("H5").Value = "'35000.00"
$ws.Range("H6").Value = "'2923.50"
$ws.Range("H7").Value = "'2182.26"
$ws.Range("H8").Value = "'961.30"
$ws.Range("H9").Value = "'12.00"
$ws.Range("H10").Value = "'4486.18"
$ws.Range("H11").Value = "'27860.98"
$ws.Range("H12").Value = "'35205.16"
$ws.Range("H13").Value = "'14405.44"
$ws.Range("H14").Value = "'9005.72"
$ws.Range("H15").Value = "'996.75"
$ws.Range("H16").Value = "'6240.00"
$ws.Range("H17").Value = "'1542.50"
$ws.Range("H18").Value = "'4061.09"
$ws.Range("H19").Value = "'80.00"
$ws.Range("H20").Value = "'490.50"
$ws.Range("H21").Value = "'655.19"
$ws.Range("H22").Value = "'6352.35"
$ws.Range("H23").Value = "'600.00"
$ws.Range("H24").Value = "'1200.00"
$ws.Range("H25").Value = "'648.00"
$ws.Range("H26").Value = "'18.00"
$ws.Range("H27").Value = "'24.20"
$ws.Range("H28").Value = "'1484.40"
$ws.Range("H29").Value = "'2040.00"
$ws.Range("H30").Value = "'114.42"
$ws.Range("H31").Value = "'1640.25"
$ws.Range("H32").Value = "'7297.48"
$ws.Range("H33").Value = "'9400.90"
$ws.Range("H34").Value = "'38.62"
$ws.Range("H35").Value = "'387.15"
$ws.Range("H36").Value = "'825.22"
$ws.Range("H37").Value = "'75.00"
$ws.Range("H38").Value = "'473.97"
$ws.Range("H39").Value = "'160.00"
$ws.Range("H40").Value = "'43.00"
$ws.Range("H41").Value = "'7043.44"
$ws.Range("H42").Value = "'815.00"
$ws.Range("H43").Value = "'861.00"
$ws.Range("H44").Value = "'488.00"
$ws.Range("H45").Value = "'516.60"
$ws.Range("H46").Value = "'1616.00"
$ws.Range("H47").Value = "'2176.00"
$ws.Range("H48").Value = "'2596.00"
$ws.Range("H49").Value = "'620.00"
$ws.Range("H50").Value = "'189.78"
$ws.Range("H51").Value = "'2200.00"
$ws.Range("H52").Value = "'6392.43"
$ws.Range("H53").Value = "'1300.00"
$ws.Range("H54").Value = "'355.00"
$ws.Range("H55").Value = "'2064.42"
$ws.Range("H56").Value = "'35.00"
$ws.Range("H57").Value = "'461.40"
$ws.Range("H58").Value = "'9270.80"
$ws.Range("H59").Value = "'4200.00"
$ws.Range("H60").Value = "'10681.60"
$ws.Range("H61").Value = "'93.00"
$ws.Range("H62").Value = "'552.30"
$ws.Range("H63").Value = "'470.00"
$ws.Range("H64").Value = "'12.00"
$ws.Range("H65").Value = "'639.00"
$ws.Range("H66").Value = "'519.00"
$ws.Range("H67").Value = "'1835.00"
$ws.Range("H68").Value = "'400.00"
$ws.Range("H69").Value = "'1450.00"
$ws.Range("H70").Value = "'110.50"
$ws.Range("H71").Value = "'460.00"
$ws.Range("H72").Value = "'2832.40"
$ws.Range("H73").Value = "'2177.90"
$ws.Range("H74").Value = "'1810.11"
$ws.Range("H75").Value = "'6753.86"
$ws.Range("H76").Value = "'74.78"
$ws.Range("H77").Value = "'177.00"
$ws.Range("H78").Value = "'1309.99"
$ws.Range("H79").Value = "'1270.00"
$ws.Range("H80").Value = "'232.70"
$ws.Range("H81").Value = "'353.00"
$ws.Range("H82").Value = "'11461.00"
$ws.Range("H83").Value = "'500.00"
$ws.Range("H84").Value = "'230.00"
$ws.Range("H85").Value = "'720.00"
$ws.Range("H86").Value = "'675.00"
$ws.Range("H87").Value = "'300.00"
$ws.Range("H88").Value = "'243.20"
$ws.Range("H89").Value = "'240.00"
$ws.Range("H90").Value = "'5304.00"
$ws.Range("H91").Value = "'200.00"
$ws.Range("H92").Value = "'750.00"
$ws.Range("H93").Value = "'220.00"
$ws.Range("H94").Value = "'2263.18"
$ws.Range("H95").Value = "'351.00"
$ws.Range("H96").Value = "'2481.83"
$ws.Range("H97").Value = "'3630.00"
$ws.Range("H98").Value = "'230.00"
$ws.Range("H99").Value = "'2130.55"
$ws.Range("H100").Value = "'11416.86"
$ws.Range("H101").Value = "'1480.00"
$ws.Range("H102").Value = "'500.00"
$ws.Range("H103").Value = "'397.80"
$ws.Range("H104").Value = "'120.00"
$ws.Range("H105").Value = "'9189.25"
$ws.Range("H106").Value = "'238.30"
$ws.Range("H107").Value = "'855.00"
$ws.Range("H108").Value = "'45.00"
$ws.Range("H109").Value = "'4300.00"
$ws.Range("H110").Value = "'3147.00"
$ws.Range("H111").Value = "'1153.00"
$ws.Range("H112").Value = "'27507.85"
$ws.Range("H113").Value = "'294.20"
$ws.Range("H114").Value = "'3831.07"
$ws.Range("H115").Value = "'40.00"
$ws.Range("H116").Value = "'20.50"
$ws.Range("H117").Value = "'4753.55"
$ws.Range("H118").Value = "'147.40"
$ws.Range("H119").Value = "'427.50"
$ws.Range("H120").Value = "'10473.50"
$ws.Range("H121").Value = "'456.00"
$ws.Range("H122").Value = "'3536.83"
$ws.Range("H123").Value = "'825.45"
$ws.Range("H124").Value = "'3375.00"
$ws.Range("H125").Value = "'2542.25"
$ws.Range("H126").Value = "'28.40"
$ws.Range("H127").Value = "'40500.00"
$ws.Range("H128").Value = "'979.26"
$ws.Range("H129").Value = "'3850.00"
$ws.Range("H130").Value = "'779.25"
$ws.Range("H131").Value = "'3100.00"
$ws.Range("H132").Value = "'1082.00"
$ws.Range("H133").Value = "'450.00"
